$wb = $excel.ActiveWorkbook

function Add-LifterRow {
    param(
        $ws,
        $row,
        $timeVal,
        $timeFmt,
        $bVal,
        $cVal,
        $dVal,
        $eVal,
        $fVal,
        $gVal,
        $gAsText,
        $hVal,
        $iVal
    )

    $ws.Cells.Item($row, 1).NumberFormat = $timeFmt
    $ws.Cells.Item($row, 1).Value = $timeVal

    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    $ws.Cells.Item($row, 6).Value = $fVal

    if ($gAsText) {
        $ws.Cells.Item($row, 7).Value = "'" + $gVal
    } else {
        $ws.Cells.Item($row, 7).Value = $gVal
    }

    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

$gNum = [double]"5.68631262647114e+23"

# Sheet 1: ROW50-FE-LIFTER -- add row 41
$ws1 = $wb.Worksheets.Item(1)
$fmt1 = $ws1.Cells.Item(40, 1).NumberFormat()
Add-LifterRow $ws1 41 45744.67248324074 $fmt1 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x6e" "0xe" 400 $gNum $false 366 14

# Sheet 2: ROW50-MID-LIFTER -- add row 43
$ws2 = $wb.Worksheets.Item(2)
$fmt2 = $ws2.Cells.Item(42, 1).NumberFormat()
$gText = "568631262647113771663628"
Add-LifterRow $ws2 43 45744.64055555555 $fmt2 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x72" "0x19" 400 $gText $true 370 25

# Sheet 3: ROW11-FE-LIFTER -- add row 41
$ws3 = $wb.Worksheets.Item(3)
$fmt3 = $ws3.Cells.Item(40, 1).NumberFormat()
Add-LifterRow $ws3 41 45744.69266969908 $fmt3 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x6e" "0x14" 400 $gNum $false 366 20

# Sheet 4: ROW11-MID-LIFTER -- add row 41
$ws4 = $wb.Worksheets.Item(4)
$fmt4 = $ws4.Cells.Item(40, 1).NumberFormat()
Add-LifterRow $ws4 41 45744.83702253472 $fmt4 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x72" "0x19" 400 $gNum $false 370 25
